$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain text (non-numeric-looking) - set directly.
$ws.Range("D2").Value = "56.110.25"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.302.83"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -3.84%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "2.315.29"
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("E10").Value = "  -3.01%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -2.48%  "
$ws.Range("E13").Value = "  -2.08%  "
$ws.Range("D14").Value = "2.713.23"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").Value = "56.041.02"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "2.317.63"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +7.20%  "
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "0.0₃0710"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  -3.13%  "

# Cells whose new values look like plain numbers to Excel's type-inference; force text
# by temporarily marking the range as Text-formatted, then clear the leftover format
# so the cell keeps its original (default) style.
$numericLookingCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D12", "D13", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "516.48"
$ws.Range("D6").Value = "130.36"
$ws.Range("D7").Value = "0.997"
$ws.Range("D8").Value = "0.529"
$ws.Range("D10").Value = "0.0992"
$ws.Range("D12").Value = "5.21"
$ws.Range("D13").Value = "0.335"
$ws.Range("D15").Value = "23.19"
$ws.Range("D19").Value = "328.40"
$ws.Range("D20").Value = "10.29"
$ws.Range("D21").Value = "4.11"
$ws.Range("D22").Value = "6.66"
$ws.Range("D23").Value = "0.999"
$ws.Range("D24").Value = "60.86"
$ws.Range("D25").Value = "0.163"
$ws.Range("D26").Value = "8.55"
$ws.Range("D27").Value = "0.991"
$ws.Range("D28").Value = "1.32"
$ws.Range("D29").Value = "167.04"
$ws.Range("D32").Value = "6.03"
$ws.Range("D33").Value = "18.17"
$ws.Range("D34").Value = "0.999"
$ws.Range("D36").Value = "1.23"
$ws.Range("D37").Value = "0.879"
$ws.Range("D38").Value = "3.86"
$ws.Range("D39").Value = "1.56"
$ws.Range("D40").Value = "38.45"
$ws.Range("D41").Value = "147.28"
$ws.Range("D42").Value = "0.373"
$ws.Range("D43").Value = "283.38"
$ws.Range("D44").Value = "3.54"
$ws.Range("D45").Value = "5.03"
$ws.Range("D46").Value = "0.0922"
$ws.Range("D47").Value = "0.0494"
$ws.Range("D48").Value = "0.553"
$ws.Range("D49").Value = "18.07"
$ws.Range("D50").Value = "0.376"
$ws.Range("D51").Value = "0.0212"
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).ClearFormats()
}
